# Updates cryptocurrency price/volume data in-place, preserving
# the existing text formatting of the Price (D) and Volume(1h) (E) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.065.60"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.263.05"
$ws.Range("E3").Value = "  -3.40%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "297.92"
$ws.Range("E5").Value = "  -2.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.10"
$ws.Range("E6").Value = "  -7.06%  "

$ws.Range("E7").Value = "  -3.04%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -4.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.06"
$ws.Range("E10").Value = "  -5.91%  "

$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.86"
$ws.Range("E12").Value = "  -8.46%  "

$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.613.85"
$ws.Range("E15").Value = "  -3.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.48"
$ws.Range("E16").Value = "  -2.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.266.44"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.774"
$ws.Range("E18").Value = "  -4.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.050.20"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("E20").Value = "  -2.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("E21").Value = "  -3.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.39"
$ws.Range("E22").Value = "  -3.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.66"
$ws.Range("E23").Value = "  -2.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.27"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("E25").Value = "  -4.96%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -4.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.79"
$ws.Range("E28").Value = "  -6.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  -1.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.19"
$ws.Range("E30").Value = "  +3.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.77"
$ws.Range("E31").Value = "  -4.30%  "

$ws.Range("E32").Value = "  -3.79%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -3.09%  "

$ws.Range("E36").Value = "  -5.80%  "

$ws.Range("E37").Value = "  -5.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.17"
$ws.Range("E38").Value = "  -7.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.77"
$ws.Range("E39").Value = "  -6.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0986"
$ws.Range("E40").Value = "  -3.98%  "

$ws.Range("E41").Value = "  -3.82%  "

$ws.Range("E42").Value = "  -9.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.43"
$ws.Range("E43").Value = "  -1.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.959.43"
$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("E45").Value = "  -2.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.29"
$ws.Range("E46").Value = "  -8.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  -7.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -5.73%  "

$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.487.39"
$ws.Range("E50").Value = "  -2.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.69"
$ws.Range("E51").Value = "  -7.74%  "
